# Daily attendance processing - 2025-12-19 18:36:57
# Rotate the "Recorded By" (column G) entries left by one position for every
# row that contains more than one comma-separated recorder so that the first
# recorder moves to the end of the list. Single-value cells are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value()

    if ($null -ne $value -and $value -is [string] -and $value.Contains(", ")) {
        $parts = $value.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        if ($parts.Length -gt 1) {
            $rotated = @($parts[1..($parts.Length - 1)]) + @($parts[0])
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
